$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the existing "sum" header (G1) onto the new
# "Save" header cell (H1) so it reuses the same header style, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values for the data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
